$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: id 4 -> 3, type truck -> car/taxi, entry_time 17:59:54 -> 02:28:19
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "car/taxi"
$ws.Range("E2").Value = "02:28:19"

# Update row 3: id 9 -> 24, entry_time 18:00:07 -> 02:28:48
$ws.Range("A3").Value = 24
$ws.Range("E3").Value = "02:28:48"

# Add new row 4
$ws.Range("A4").Value = 16
$ws.Range("B4").Value = "car/taxi"
$ws.Range("C4").Value = "west"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "02:29:29"

# Add new row 5
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "car/taxi"
$ws.Range("C5").Value = "south"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "02:29:33"

# Add new row 6
$ws.Range("A6").Value = 9
$ws.Range("B6").Value = "car/taxi"
$ws.Range("C6").Value = "west"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "02:29:42"

# Add new row 7
$ws.Range("A7").Value = 33
$ws.Range("B7").Value = "car/taxi"
$ws.Range("C7").Value = "east"
$ws.Range("D7").Value = "north"
$ws.Range("E7").Value = "02:29:45"
